$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30: Binary Tree Inorder Traversal (Morris Traversal) ---
$ws.Range("A30").Value = 94
$ws.Range("B30").Value = "LC"
$ws.Range("C30").Value = "Binary Tree Inorder Traversal(Morris Traversal)"
$rng30 = $ws.Range("C30").Characters(31, 16)
$rng30.Font.Bold = $true
$ws.Range("D30").Value = "Java/Python"
$ws.Range("E30").Value = "Medium"

# --- Row 31: Binary Tree Preorder Traversal-(Morris Traversal) ---
$ws.Range("A31").Value = 144
$ws.Range("B31").Value = "LC"
$ws.Range("C31").Value = "Binary Tree Preorder Traversal-(Morris Traversal)"
$ws.Range("D31").Value = "Java/Python"
$ws.Range("E31").Value = "Medium"

# --- Row 32: Flatten Binary Tree to Linked List ---
$ws.Range("A32").Value = 114
$ws.Range("B32").Value = "LC"
$ws.Range("C32").Value = "Flatten Binary Tree to Linked List"
$ws.Range("D32").Value = "Java"
$ws.Range("E32").Value = "Medium"

# Match the author's final view / selection state
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("E32").Select()
